# "change worker menu, and get_sales_report"
# Adds two new sales rows (a skirts sale and a blouses sale/refund) to the
# "Sales01" worksheet, as produced by the updated get_sales_report routine.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sales01")

# Row 10: 2019-01-04, product code 12, "skirts", qty 1, price 80, receipt 4
$ws.Cells.Item(10, 1).Value = 2019
$ws.Cells.Item(10, 2).Value = 1
$ws.Cells.Item(10, 3).Value = 4
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = "skirts"
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 80
$ws.Cells.Item(10, 8).Value = 4

# Row 11: 2019-01-04, product code 2, "blouses", qty 2, price 119.9, receipt 4
$ws.Cells.Item(11, 1).Value = 2019
$ws.Cells.Item(11, 2).Value = 1
$ws.Cells.Item(11, 3).Value = 4
$ws.Cells.Item(11, 4).Value = 2
$ws.Cells.Item(11, 5).Value = "blouses"
$ws.Cells.Item(11, 6).Value = 2
$ws.Cells.Item(11, 7).Value = 119.9
$ws.Cells.Item(11, 8).Value = 4
